$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 30   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Crime-complaint table updates (rows 15-29) ---
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = "'0"
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = -50
$ws.Range("N15").Value = -85.714285714285
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 150
$ws.Range("I16").Value = 19
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = 72.727272727272
$ws.Range("L16").Value = 90
$ws.Range("M16").Value = -47.222222222222
$ws.Range("N16").Value = -86.428571428571
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 10.526315789473
$ws.Range("I17").Value = 26
$ws.Range("J17").Value = 28
$ws.Range("K17").Value = -7.142857142857
$ws.Range("L17").Value = 4
$ws.Range("M17").Value = -21.212121212121
$ws.Range("N17").Value = -65.78947368421
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 45.454545454545
$ws.Range("I18").Value = 21
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = 61.538461538461
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = 5
$ws.Range("N18").Value = -66.129032258064
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 3.846153846153
$ws.Range("I19").Value = 40
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 14.285714285714
$ws.Range("L19").Value = 48.148148148148
$ws.Range("M19").Value = 81.818181818181
$ws.Range("N19").Value = 25
$ws.Range("C20").Value = "'0"
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = -57.142857142857
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -47.058823529411
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -87.323943661971
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -13.636363636363
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 8.974358974358
$ws.Range("I21").Value = 116
$ws.Range("J21").Value = 107
$ws.Range("K21").Value = 8.411214953271
$ws.Range("L21").Value = 30.337078651685
$ws.Range("M21").Value = -6.451612903225
$ws.Range("N21").Value = -70.483460559796
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 10
$ws.Range("J23").Value = 7
$ws.Range("K23").Value = 42.857142857142
$ws.Range("L23").Value = 25
$ws.Range("M23").Value = 11.111111111111
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 1.724137931034
$ws.Range("I24").Value = 77
$ws.Range("J24").Value = 66
$ws.Range("K24").Value = 16.666666666666
$ws.Range("L24").Value = 50.980392156862
$ws.Range("M24").Value = 18.461538461538
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 12.5
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 3.333333333333
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 37
$ws.Range("K25").Value = 5.405405405405
$ws.Range("L25").Value = 69.565217391304
$ws.Range("M25").Value = -49.350649350649
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = "'0"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = -100
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = -33.333333333333
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = -75
$ws.Range("L27").Value = -75
$ws.Range("M28").Value = -71.428571428571
$ws.Range("N28").Value = -88.888888888888
$ws.Range("M29").Value = -71.428571428571
$ws.Range("N29").Value = -87.5
